# Updated cryptos list (price + 1h volume-change refresh), plus a rank swap
# between Monero and ImmutableX (rows 36/37).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.589.31"
$ws.Range("E2").Value = "  -3.26%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.346.31"
$ws.Range("E3").Value = "  -2.76%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "566.55"
$ws.Range("E5").Value = "  -2.35%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.73"
$ws.Range("E6").Value = "  -0.75%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.484"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.90"
$ws.Range("E9").Value = "  -1.00%  "
$ws.Range("E10").Value = "  -1.11%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.416"
$ws.Range("E11").Value = "  +1.69%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "3.912.26"
$ws.Range("E12").Value = "  -2.93%  "
$ws.Range("E13").Value = "  +1.20%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.70"
$ws.Range("E14").Value = "  -2.00%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.347.20"
$ws.Range("E15").Value = "  -2.69%  "
$ws.Range("E16").Value = "  -1.24%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "60.592.40"
$ws.Range("E17").Value = "  -3.39%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.29"
$ws.Range("E18").Value = "  -1.14%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.52"
$ws.Range("E19").Value = "  -1.09%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.92"
$ws.Range("E20").Value = "  -1.57%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "376.58"
$ws.Range("E21").Value = "  -2.60%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.560"
$ws.Range("E22").Value = "  -0.38%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "74.73"
$ws.Range("E23").Value = "  -0.76%  "
$ws.Range("E24").Value = "  -0.14%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.494.51"
$ws.Range("E25").Value = "  -2.43%  "
$ws.Range("E26").Value = "  -4.93%  "
$ws.Range("E27").Value = "  -3.94%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.994"
$ws.Range("E28").Value = "  -0.57%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.34"
$ws.Range("E29").Value = "  -3.96%  "
$ws.Range("E30").Value = "  -1.11%  "
$ws.Range("E31").Value = "  +0.02%  "
$ws.Range("E32").Value = "  -3.37%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "22.92"
$ws.Range("E33").Value = "  -1.21%  "
$ws.Range("E34").Value = "  -3.77%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.32"
$ws.Range("E35").Value = "  -0.63%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.83"
$ws.Range("E38").Value = "  -1.61%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "27.90"
$ws.Range("E39").Value = "  -12.79%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.378.36"
$ws.Range("E40").Value = "  -2.88%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0749"
$ws.Range("E41").Value = "  -3.14%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.757"
$ws.Range("E42").Value = "  -3.55%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.30"
$ws.Range("E43").Value = "  -1.28%  "
$ws.Range("E44").Value = "  -5.14%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.14"
$ws.Range("E45").Value = "  -3.51%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.459.08"
$ws.Range("E46").Value = "  -4.06%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.67"
$ws.Range("E47").Value = "  -3.54%  "
$ws.Range("E48").Value = "  -0.01%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "22.31"
$ws.Range("E49").Value = "  -1.50%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0261"
$ws.Range("E50").Value = "  -1.84%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.819"
$ws.Range("E51").Value = "  +0.67%  "

# Rows 36/37 swap rank order: Monero drops below ImmutableX, each refreshed
# with its own new price / 1h volume-change.
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.55"
$ws.Range("E36").Value = "  -4.75%  "

$ws.Range("B37").Value = "Monero"
$ws.Range("C37").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "168.24"
$ws.Range("E37").Value = "  -0.74%  "
